# Add a merged "Note" row under the customer table, bold + centered,
# matching the sampleCustomers.xlsx "ok pre final code" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string row: Sr./name/email/cell/address columns stay blank,
# only A11 carries the note text (cells get merged across A:D below).
$note = $ws.Range("A11")
$note.Value = "Note:All fields are required for every item."

# Bold + centered formatting for the note (format the anchor cell before
# merging so every cell in the merged block picks up the same style).
$note.Font.Bold = $true
$note.HorizontalAlignment = -4108  # xlCenter

# Merge A11:D11 into a single cell spanning the four data columns.
$ws.Range("A11:D11").Merge()

# Match the workbook's new selection (the merged note cell/range).
$ws.Range("A11:D11").Select() | Out-Null

# Page setup was switched to explicit portrait orientation.
$ws.PageSetup.Orientation = 1  # xlPortrait
